# Resume formatting update:
#  - widen the indent/date/body columns of the main layout table
#  - add right padding + grey color to the date column on EXPERIENCE rows
#  - use a non-breaking space before the en-dash in date ranges so they
#    don't wrap awkwardly

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# --- Column widths -------------------------------------------------------
# New widths (in twentieths of a point / dxa -> pt by dividing by 20):
#   col1 (indent)      28  -> 144
#   col2 (label/date) 1584 -> 1944
#   col3 (body)       8683 -> 8208
# The header row's second cell spans columns 2+3 (gridSpan=2), so its own
# width becomes 1944+8208 = 10152. Set that merged cell's width first so
# that the later per-row writes for the un-merged rows are the ones that
# finally win for the shared table grid / column widths.
$t.Cell(1,1).Width = 144 / 20
$t.Cell(1,2).Width = 10152 / 20

for ($r = 2; $r -le $t.Rows.Count; $r++) {
    $row = $t.Rows.Item($r)
    $row.Cells.Item(1).Width = 144 / 20
    $row.Cells.Item(2).Width = 1944 / 20
    $row.Cells.Item(3).Width = 8208 / 20
}

# --- Date column formatting (EXPERIENCE rows) -----------------------------
$nbsp = [char]0x00A0
$endash = [char]0x2013
$searchStr = " $endash"
$replaceStr = "$nbsp$endash"

$dateRows = 4, 5, 6, 7, 8
foreach ($r in $dateRows) {
    $cell = $t.Cell($r, 2)

    # Add right padding (tcMar/right = 216 dxa = 10.8pt)
    $cell.RightPadding = 216 / 20

    # Replace the regular space before the en-dash with a non-breaking space
    $cell.Range.Find.Execute($searchStr, $true, $false, $false, $false, $false, $true, 1, $false, $replaceStr, 2) | Out-Null

    # Apply the grey color to the date text run (exclude the trailing
    # paragraph/cell-end marks from the range)
    $rng = $cell.Range
    $textRng = $d.Range($rng.Start, $rng.End - 1)
    $textRng.Font.Color = 6710886
}
